$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.111.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.02%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.888.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.65%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'307.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.12%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.5151"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.86%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3722"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.15%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07209"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.67%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9030"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.44%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'21.05"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.89%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07629"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.85%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.903.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.42%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'94.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.41%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +0.88%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -0.04%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.000008512"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.10%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D19").Value = "'0.9996"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.00%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'27.162.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.04%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.051"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.72%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'2.141.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.20%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'6.422"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.07%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'145.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.95%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.792"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.20%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'18.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.32%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.176"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +5.90%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'114.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.53%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.983"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +6.95%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.826"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +4.25%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.09209"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.11%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.05065"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.33%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +4.68%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.7629"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.78%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.003"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.75%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.277"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.11%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.575"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.01%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.5612"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +5.39%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.18%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.18%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'9.023"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +7.85%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'118.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.46%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'6.578"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.76%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.1507"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.49%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4802"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +3.50%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'10.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.04%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.9998"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.05%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.580"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.40%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'37.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.24%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'63.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.33%  "
$ws.Range("E51").Style = "Normal"
